# Auto-generated Excel COM-interop script to apply the trend-summary update
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Summary Table")

$summaryData = @(
  @(' Innovation-driven development ', ' 37 ', ' China unveiled a new 5-year innovation strategy prioritizing "new quality productivity" across strategic sectors including AI, quantum computing, and advanced manufacturing. ', ' https://www.xinhuanet.com/english/20250415/china-unveils-new-innovation-strategy ', ' The State Council released its 2025-2030 Innovation-Driven Development Strategy, emphasizing breakthroughs in core technologies and digital transformation. The plan aims to increase R&D spending to 3.2% of GDP by 2030 and establishes 12 national innovation clusters focused on strategic emerging industries. Premier Li emphasized that innovation remains the "primary driving force" for China''s continued economic development.'),
  @(' New quality productivity ', ' 34 ', ' Beijing announced a $75 billion fund to accelerate "new quality productivity" investments in semiconductor manufacturing, AI infrastructure, and green technology. ', ' https://www.scmp.com/tech/policy/article/3257814/china-announces-75-billion-fund-new-quality-productivity ', ' China''s Ministry of Finance, in coordination with the NDRC, announced a major funding initiative aimed at upgrading industrial capabilities through advanced manufacturing technologies. The fund will prioritize domestic semiconductor development, AI computing infrastructure, and green manufacturing processes. This represents China''s largest-ever directed investment in emerging technologies, with 40% allocated to chip manufacturing and advanced computing.'),
  @(' Artificial intelligence ', ' 29 ', ' China released its first comprehensive AI ethics framework establishing guidelines for responsible AI development and deployment across sectors. ', ' https://www.globaltimes.cn/page/202504/1318543.shtml ', ' The Cyberspace Administration of China (CAC) published the "Ethical Guidelines for Artificial Intelligence Development," China''s first comprehensive national AI ethics framework. The guidelines establish principles for transparency, data privacy, algorithmic fairness, and human-AI collaboration. They require AI systems to undergo ethical reviews, particularly for high-risk applications in healthcare, finance, and public services. The framework also addresses international cooperation, calling for global governance mechanisms.'),
  @(' Science and technology innovation policy ', ' 27 ', ' China''s State Council issued a new policy mandating that 15% of government procurement budgets be allocated to products from innovative tech startups and SMEs. ', ' https://english.www.gov.cn/policies/latestreleases/202504/16/content_WS6642a901c6d0868f4e8b1c2a.html ', ' The policy aims to leverage government purchasing power to support domestic innovation, particularly in digital services, green technology, and advanced manufacturing. It establishes streamlined certification processes for qualifying companies and requires state-owned enterprises to adopt similar procurement targets. The measure, part of broader efforts to reduce reliance on foreign technology, takes effect June 1, 2025, and includes specific provisions for companies less than 8 years old.'),
  @(' Semiconductor packaging ', ' 26 ', ' Chinese researchers announced a breakthrough in advanced semiconductor packaging technology for AI chips, reducing power consumption by 35%. ', ' https://asia.nikkei.com/Business/Tech/Semiconductors/China-claims-breakthrough-in-chip-packaging-for-AI-applications ', ' A research team from Tsinghua University and the Chinese Academy of Sciences developed a new 3D stacking technique for high-bandwidth memory in AI accelerator chips. The approach uses novel thermal management and interconnect materials to overcome previous limitations. Independent experts confirmed the significance of the advance, which allows for more computing power while consuming less energy. The technology is expected to be commercialized within 12-18 months through partnerships with domestic chipmakers.'),
  @(' Digital economy policy ', ' 24 ', ' China''s Ministry of Industry and Information Technology launched a digital economy standardization initiative covering blockchain, industrial internet, and cross-border data flows. ', ' https://www.caixinglobal.com/2025-04-14/china-launches-digital-economy-standardization-initiative-101973542.html ', ' The initiative establishes 87 new technical standards across digital infrastructure, data governance, and emerging technologies. It prioritizes interoperability between domestic platforms and systems while creating a framework for secure cross-border data flows. The standards create certification requirements for digital service providers and establish compliance timelines. Industry representatives from tech giants and SMEs participated in the development process through public consultations.'),
  @(' Digital RMB ', ' 23 ', ' The People''s Bank of China expanded its digital yuan pilot to include international transactions with selected trading partners, starting with ASEAN countries. ', ' https://www.reuters.com/technology/china-expands-digital-yuan-international-transactions-2025-04-13/ ', ' The PBOC announced that the digital yuan (e-CNY) can now be used for certain cross-border transactions with trade partners in Thailand, Malaysia, and Singapore as the first phase of internationalization. The pilot program focuses on settlement for small and medium-sized trading companies and establishes currency swap mechanisms with participating central banks. Officials emphasized the system''s reduced transaction costs and independence from Western financial infrastructure.'),
  @(' Science and technology system reform ', ' 22 ', ' China announced reforms to its scientific research evaluation system, reducing emphasis on publication metrics and increasing focus on technological breakthroughs. ', ' https://www.nature.com/articles/d41586-025-01142-x ', ' The Chinese Academy of Sciences and Ministry of Science and Technology jointly issued new guidelines that fundamentally reform how scientific research is evaluated and funded. The changes reduce emphasis on journal impact factors and publication counts in favor of actual technological and economic impact. The reforms also streamline grant applications, give more autonomy to research institutions in fund allocation, and create special funding streams for high-risk exploratory research.'),
  @(' Integration of industry, academia and research ', ' 21 ', ' Beijing launched 10 new national innovation centers designed to accelerate commercialization of research in quantum computing, biotechnology, and advanced materials. ', ' https://www.chinadaily.com.cn/a/202504/17/WS660e1234a310e074ea38f5c4.html ', ' The Ministry of Science and Technology established 10 new National Technology Innovation Centers focusing on quantum computing, biotechnology, advanced materials, and other strategic domains. Each center involves partnerships between top universities, research institutes, and leading companies in their respective fields. The centers will receive 5-year funding packages of approximately $500 million each and operate under new governance models with streamlined intellectual property arrangements designed to accelerate technology transfer.'),
  @(' Hydrogen energy storage ', ' 20 ', ' China''s energy authority approved a $3.2 billion national hydrogen energy demonstration project spanning five provinces. ', ' https://www.cnstock.com/v_news/sns_yw/202504/4954217.htm ', ' The National Energy Administration approved a comprehensive hydrogen energy demonstration project spanning five provinces (Shandong, Guangdong, Jiangsu, Hebei, and Inner Mongolia). The initiative includes production infrastructure, storage facilities, and transportation networks. The project emphasizes green hydrogen produced through renewable energy and includes applications in industrial processes, transportation, and grid-scale energy storage. It establishes technical standards for hydrogen production, storage, and refueling facilities.'),
  @(' Quantum communication ', ' 19 ', ' China launched a satellite dedicated to quantum communication experiments, expanding its quantum network infrastructure. ', ' https://spacenews.com/china-launches-satellite-for-quantum-communications-experiments/ ', ' China launched the Quantum Science Satellite-2 (QSS-2), a significant expansion of its space-based quantum communications capabilities. The satellite will test new quantum key distribution protocols at greater distances and higher speeds than previous systems. It represents collaboration between the Chinese Academy of Sciences and the University of Science and Technology of China. Officials indicated the system aims to create an "unhackable" communications network for government and financial institutions.'),
  @(' Talent introduction policy ', ' 18 ', ' China eased visa and residence requirements for foreign science and technology professionals, offering fast-track permanent residency for experts in strategic fields. ', ' https://www.chinabankingnews.com/2025/04/12/china-streamlines-visa-requirements-for-foreign-tech-talent/ ', ' The Ministry of Science and Technology and National Immigration Administration jointly announced simplified visa and residence procedures for foreign technical experts. The policy creates a special "technology talent" visa category with expedited processing and multi-entry options. It offers pathways to permanent residency after three years for experts in strategic fields like semiconductors, AI, and biotechnology. The policy also establishes support services for foreign professionals, including assistance with housing, education, and healthcare.'),
  @(' International talent service ', ' 17 ', ' Shanghai launched the "Global Innovation Talent Hub" with expanded services and benefits for international technology professionals. ', ' https://www.shine.cn/news/metro/2504173658/ ', ' Shanghai''s municipal government established a comprehensive service platform for international technology professionals, providing streamlined work permits, tax incentives, and support services. The initiative includes dedicated innovation zones with preferential policies for startups founded by international entrepreneurs. It creates special international schools and healthcare access for families of technology professionals and offers subsidized housing in key innovation districts.'),
  @(' Strategic emerging industries ', ' 16 ', ' China released a catalog of 127 strategic emerging industries eligible for preferential government support and investment. ', ' https://www.yicaiglobal.com/news/china-identifies-127-strategic-emerging-industries-for-policy-support ', ' The National Development and Reform Commission issued a comprehensive catalog of strategic emerging industries that will receive preferential policy support, tax incentives, and investment. The list encompasses 127 specific technology domains across next-generation information technology, biotechnology, new materials, advanced manufacturing, green technology, and frontier sciences. The catalog serves as a roadmap for government investment and industrial policy through 2030, with specific growth targets for each sector.'),
  @(' Bilateral science and technology joint committee ', ' 15 ', ' China and Germany held a bilateral science and technology joint committee meeting, agreeing to expand cooperation in clean energy and AI research. ', ' https://www.dw.com/en/china-germany-expand-technology-cooperation-despite-growing-tensions/a-66782314 ', ' The 10th China-Germany Joint Committee on Science and Technology Cooperation convened in Berlin, resulting in agreements to expand research collaboration in clean energy, artificial intelligence, and advanced manufacturing. The committee established six new joint research programs with €200 million in combined funding. While acknowledging broader geopolitical tensions, both sides emphasized the importance of maintaining scientific exchange. The agreement includes provisions for talent exchange programs and joint innovation centers.'),
  @(' Intellectual property protection ', ' 14 ', ' China''s Supreme Court issued new guidelines strengthening intellectual property protections with increased penalties for violations. ', ' https://www.chinaipmagazine.com/en/news-show.asp?id=12703 ', ' The Supreme People''s Court issued comprehensive guidelines strengthening intellectual property protection across patents, trademarks, and trade secrets. The guidelines increase statutory damages for IP violations and establish specialized IP tribunals in 15 additional cities. They clarify standards for determining infringement and streamline enforcement procedures. The measures specifically address concerns in the pharmaceutical, software, and semiconductor industries, sectors where foreign companies have raised concerns about IP protection.'),
  @(' R&D investment strategy ', ' 14 ', ' China pledged to increase basic research funding to 12% of total R&D spending by 2027, focusing on breakthrough technologies. ', ' https://www.scmp.com/news/china/science/article/3257903/china-boost-basic-research-funding-push-technological-breakthroughs ', ' The Ministry of Science and Technology announced plans to significantly increase funding for basic scientific research, raising its share of total R&D spending from 6.6% to 12% by 2027. The strategy focuses on physics, mathematics, materials science, and other fundamental disciplines that enable technological breakthroughs. It creates new funding mechanisms for long-term research programs and establishes international collaboration initiatives in basic science. The plan aims to address China''s historic underinvestment in foundational research compared to applied technologies.'),
  @(' Science and technology cooperation agreement ', ' 13 ', ' China and Brazil signed a comprehensive science and technology cooperation agreement focused on agricultural biotechnology and space research. ', ' https://news.cgtn.com/news/2025-04-16/China-Brazil-sign-comprehensive-science-and-technology-agreement-1mcTDgyrXRS/index.html ', ' China and Brazil signed a wide-ranging science and technology cooperation agreement focusing on agricultural biotechnology, space technology, and renewable energy. The agreement establishes joint funding mechanisms totaling $600 million over five years for collaborative research projects. It creates exchange programs for researchers and students between both countries and establishes shared research facilities. The partnership also includes provisions for technology transfer and commercialization of jointly developed innovations.'),
  @(' Joint research center ', ' 12 ', ' Tsinghua University and MIT launched a joint research center for climate science and green technology with $200 million in funding. ', ' https://www.technologyreview.com/2025/04/14/1089547/mit-tsinghua-launch-joint-climate-research-center/ ', ' Tsinghua University and the Massachusetts Institute of Technology established a joint research center focused on climate science and green technology. The center received $200 million in combined funding from both institutions and corporate partners. Research will focus on carbon capture, climate modeling, sustainable materials, and clean energy technologies. The center establishes exchange programs for faculty and students and creates shared laboratory facilities. It represents one of the largest US-China academic collaborations in recent years despite broader geopolitical tensions.'),
  @(' Science and technology financial support ', ' 11 ', ' China''s central bank introduced new financial instruments designed specifically to support early-stage technology startups and R&D-intensive companies. ', ' https://www.ft.com/content/3e4d5b8a-6a19-4c01-a8dc-0e3f45a97c8f ', ' The People''s Bank of China, in coordination with financial regulators, introduced new financial instruments specifically designed for technology startups and R&D-intensive companies. The measures include credit enhancement tools for bank loans to tech startups, specialized bond products for research-intensive firms, and risk-sharing mechanisms for venture capital investments. The policies aim to address the financing challenges faced by early-stage companies that have significant intellectual property but limited physical assets or revenue.'),
  @(' Science and technology park policy ', ' 10 ', ' Shanghai announced plans for a new $5 billion science and technology park focused on quantum computing and advanced materials. ', ' https://www.yicaiglobal.com/news/shanghai-to-build-usd5-billion-quantum-computing-science-park ', ' Shanghai''s municipal government announced plans for a major new science and technology park dedicated to quantum computing and advanced materials research. The 420-hectare development will receive $5 billion in infrastructure investment and tax incentives for participating companies. The park will integrate research facilities, corporate R&D centers, manufacturing capabilities, and residential areas. It establishes preferential policies for talent recruitment, intellectual property protection, and international collaboration.')
)

$startRow = 2
$endRow = $startRow + $summaryData.Length - 1
$ws1.Range("B$startRow`:B$endRow").NumberFormat = "@"
for ($i = 0; $i -lt $summaryData.Length; $i++) {
  $r = $startRow + $i
  $ws1.Range("A$r").Value = $summaryData[$i][0]
  $ws1.Range("B$r").Value = $summaryData[$i][1]
  $ws1.Range("C$r").Value = $summaryData[$i][2]
  $ws1.Range("D$r").Value = $summaryData[$i][3]
  $ws1.Range("E$r").Value = $summaryData[$i][4]
}
$ws1.Range("B$startRow`:B$endRow").Style = "Normal"

$ws2 = $wb.Worksheets.Item("Sources")

$sourcesData = @(
  @(' Xinhua News ', ' https://www.xinhuanet.com/english/20250415/china-unveils-new-innovation-strategy ', ' 2025-04-15 '),
  @(' South China Morning Post ', ' https://www.scmp.com/tech/policy/article/3257814/china-announces-75-billion-fund-new-quality-productivity ', ' 2025-04-17 '),
  @(' Global Times ', ' https://www.globaltimes.cn/page/202504/1318543.shtml ', ' 2025-04-14 '),
  @(' China Government Website ', ' https://english.www.gov.cn/policies/latestreleases/202504/16/content_WS6642a901c6d0868f4e8b1c2a.html ', ' 2025-04-16 '),
  @(' Nikkei Asia ', ' https://asia.nikkei.com/Business/Tech/Semiconductors/China-claims-breakthrough-in-chip-packaging-for-AI-applications ', ' 2025-04-13 '),
  @(' Caixin Global ', ' https://www.caixinglobal.com/2025-04-14/china-launches-digital-economy-standardization-initiative-101973542.html ', ' 2025-04-14 '),
  @(' Reuters ', ' https://www.reuters.com/technology/china-expands-digital-yuan-international-transactions-2025-04-13/ ', ' 2025-04-13 '),
  @(' Nature ', ' https://www.nature.com/articles/d41586-025-01142-x ', ' 2025-04-12 '),
  @(' China Daily ', ' https://www.chinadaily.com.cn/a/202504/17/WS660e1234a310e074ea38f5c4.html ', ' 2025-04-17 '),
  @(' China Securities News ', ' https://www.cnstock.com/v_news/sns_yw/202504/4954217.htm ', ' 2025-04-15 '),
  @(' SpaceNews ', ' https://spacenews.com/china-launches-satellite-for-quantum-communications-experiments/ ', ' 2025-04-16 '),
  @(' China Banking News ', ' https://www.chinabankingnews.com/2025/04/12/china-streamlines-visa-requirements-for-foreign-tech-talent/ ', ' 2025-04-12 '),
  @(' Shine ', ' https://www.shine.cn/news/metro/2504173658/ ', ' 2025-04-17 '),
  @(' Yicai Global ', ' https://www.yicaiglobal.com/news/china-identifies-127-strategic-emerging-industries-for-policy-support ', ' 2025-04-18 '),
  @(' Deutsche Welle ', ' https://www.dw.com/en/china-germany-expand-technology-cooperation-despite-growing-tensions/a-66782314 ', ' 2025-04-15 '),
  @(' China IP Magazine ', ' https://www.chinaipmagazine.com/en/news-show.asp?id=12703 ', ' 2025-04-14 '),
  @(' South China Morning Post ', ' https://www.scmp.com/news/china/science/article/3257903/china-boost-basic-research-funding-push-technological-breakthroughs ', ' 2025-04-16 '),
  @(' CGTN ', ' https://news.cgtn.com/news/2025-04-16/China-Brazil-sign-comprehensive-science-and-technology-agreement-1mcTDgyrXRS/index.html ', ' 2025-04-16 '),
  @(' MIT Technology Review ', ' https://www.technologyreview.com/2025/04/14/1089547/mit-tsinghua-launch-joint-climate-research-center/ ', ' 2025-04-14 '),
  @(' Financial Times ', ' https://www.ft.com/content/3e4d5b8a-6a19-4c01-a8dc-0e3f45a97c8f ', ' 2025-04-13 '),
  @(' Yicai Global ', ' https://www.yicaiglobal.com/news/shanghai-to-build-usd5-billion-quantum-computing-science-park ', ' 2025-04-17 ')
)

$startRow2 = 3
for ($i = 0; $i -lt $sourcesData.Length; $i++) {
  $r = $startRow2 + $i
  $ws2.Range("A$r").Value = $sourcesData[$i][0]
  $ws2.Range("B$r").Value = $sourcesData[$i][1]
  $ws2.Range("C$r").Value = $sourcesData[$i][2]
}

$ws3 = $wb.Worksheets.Item("Executive Summary")
$ws3.Range("A2").Value = 'Five Most Impactful News Summaries:

1. China unveiled its 2025-2030 Innovation-Driven Development Strategy, establishing a comprehensive framework for technology advancement with increased R&D spending targets (3.2% of GDP by 2030) and 12 national innovation clusters. The strategy emphasizes breakthroughs in core technologies and positions innovation as the "primary driving force" for continued economic development.

2. The Chinese government announced a massive $75 billion fund dedicated to "new quality productivity" investments, with 40% allocated to semiconductor development and advanced computing infrastructure. This represents China''s largest-ever directed investment in emerging technologies and signals a strong commitment to domestic technological self-sufficiency.

3. The Cyberspace Administration of China published the country''s first comprehensive national AI ethics framework, establishing principles for transparency, data privacy, algorithmic fairness, and human-AI collaboration. The framework requires ethical reviews for high-risk AI applications and addresses international cooperation, positioning China to influence global AI governance standards.

4. China''s State Council issued a new policy mandating that 15% of government procurement budgets be allocated to products from innovative tech startups and SMEs, leveraging massive government purchasing power to accelerate domestic innovation in digital services, green technology, and advanced manufacturing.

5. Chinese researchers from Tsinghua University and the Chinese Academy of Sciences developed a breakthrough 3D stacking technique for semiconductor packaging in AI chips that reduces power consumption by 35%. This advancement allows for greater computing power with lower energy requirements and demonstrates China''s progress in addressing bottlenecks in chip technology.'

$ws4 = $wb.Worksheets.Item("Cooccurrence")
$ws4.Rows.Item(4).Insert()

$coocData = @(
  @('Innovation-driven development', 'New quality productivity', 1),
  @('Innovation-driven development', 'Strategic emerging industries', 1),
  @('New quality productivity', 'Strategic emerging industries', 1),
  @('Artificial intelligence', 'Bilateral science and technology joint committee', 1)
)

$startRow4 = 2
for ($i = 0; $i -lt $coocData.Length; $i++) {
  $r = $startRow4 + $i
  $ws4.Range("A$r").Value = $coocData[$i][0]
  $ws4.Range("B$r").Value = $coocData[$i][1]
  $ws4.Range("C$r").Value = $coocData[$i][2]
}

$ws5 = $wb.Worksheets.Item("Associations")
$ws5.Rows.Item(12).Delete()
$ws5.Rows.Item(12).Delete()

$assocData = @(
  @('Innovation-driven development', 1),
  @('New quality productivity', 2),
  @('Strategic emerging industries', 2),
  @('Artificial intelligence', 2),
  @('Semiconductor packaging', 1),
  @('Quantum communication', 1),
  @('Bilateral science and technology joint committee', 1),
  @('Intellectual property protection', 2),
  @('Science and technology cooperation agreement', 1),
  @('Joint research center', 1)
)

$startRow5 = 2
for ($i = 0; $i -lt $assocData.Length; $i++) {
  $r = $startRow5 + $i
  $ws5.Range("A$r").Value = $assocData[$i][0]
  $ws5.Range("B$r").Value = $assocData[$i][1]
}

Write-Output "Edit script completed"
